$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 213, shifting rows 213:301 down to 214:302.
$ws.Rows.Item(213).Insert()

# Fill in the new row 213 with the new data (the rest of the columns, A,B,C,E,F,G,H,I,N,O,Q,R,
# are constant across this block so they need to be copied too since Insert() leaves row 213 blank).
$ws.Range("A213").Value = 6
$ws.Range("B213").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C213").Value = "Metropolitana"
$ws.Range("D213").Value = 44489
$ws.Range("E213").Value = 13
$ws.Range("F213").Value = 100112039
$ws.Range("G213").Value = "Ciboulette"
$ws.Range("H213").Value = "Sin especificar"
$ws.Range("I213").Value = "Primera"
$ws.Range("J213").Value = 710
$ws.Range("K213").Value = 800
$ws.Range("L213").Value = 900
$ws.Range("M213").Value = 838
$ws.Range("N213").Value = "$/docena de atados"
$ws.Range("O213").Value = "Región Metropolitana"
$ws.Range("P213").Value = 279
$ws.Range("Q213").Value = 3
$ws.Range("R213").Value = "Hortaliza"
